# Update "想去人数" (interested headcount) figures in the three sheets that
# share event rows: 展览 (Exhibition), 演出 (Performance), 全部类型 (All types).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 535
$wsExhibition.Range("F5").Value = 259
$wsExhibition.Range("F6").Value = 388
$wsExhibition.Range("F8").Value = 2291
$wsExhibition.Range("F10").Value = 5695
$wsExhibition.Range("F11").Value = 138
$wsExhibition.Range("F12").Value = 373

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F3").Value = 11

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 535
$wsAllTypes.Range("F6").Value = 259
$wsAllTypes.Range("F7").Value = 388
$wsAllTypes.Range("F9").Value = 11
$wsAllTypes.Range("F11").Value = 2291
$wsAllTypes.Range("F13").Value = 5695
$wsAllTypes.Range("F14").Value = 138
$wsAllTypes.Range("F15").Value = 373
